$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("VENTA MENSUAL")

# New column G, same width look as original (width=17 in XML units)
$ws.Columns.Item(7).ColumnWidth = 16.14

# G1 header: copy format from F1 (bold header style) then set text
$ws.Cells.Item(1, 6).Copy()
$ws.Cells.Item(1, 7).PasteSpecial(-4122)
$ws.Cells.Item(1, 7).Value = "PRESUPUESTO"

# G2:G6 data cells: copy format from same-row F cell, then set value 0
for ($r = 2; $r -le 6; $r++) {
    $ws.Cells.Item($r, 6).Copy()
    $ws.Cells.Item($r, 7).PasteSpecial(-4122)
    $ws.Cells.Item($r, 7).Value = 0
}

# G7 totals row: copy format from F7, then set value 0
$ws.Cells.Item(7, 6).Copy()
$ws.Cells.Item(7, 7).PasteSpecial(-4122)
$ws.Cells.Item(7, 7).Value = 0

$excel.CutCopyMode = $false
